$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: fix "OU" -> "OO" typo and apply the yellow highlight used by the
# rest of the header row (I3:L3) so the whole legend-header row is consistent.
$ws.Range("I3").Value2 = "OO"
$ws.Range("I3:L3").Interior.Color = 65535

# --- Row 7: apply the same yellow highlight to the pin-name cells so they
# match the style already used on A7.
$ws.Range("B7:D7").Interior.Color = 65535
$ws.Range("F7:G7").Interior.Color = 65535
$ws.Range("I7").Interior.Color = 65535

# H7 and J7 already carry a (borderless) "applyBorder" format - keep that
# flag set while adding the highlight fill so they pick up their own
# (slightly different) style, matching the rest of the legend markers.
foreach ($addr in "H7", "J7") {
    $cell = $ws.Range($addr)
    $cell.Interior.Color = 65535
    $cell.Borders.Item(5).LineStyle = -4142
}

# The "RPT"/"RPR" shortcuts for the Raspberry Pi Tx/Rx pins are no longer
# called out in this row - clear their labels but keep the cell formatting.
$ws.Range("L7:M7").ClearContents()

# --- Remove the now-unused "Raspberry Pi Transmit/Receive" legend entries
# at the bottom of the sheet entirely.
$ws.Rows("20:21").Delete()

# Reset the active selection to the top-right corner of the table, like a
# user would after finishing the edits.
$ws.Range("O3").Select()
